# Auto-generated Excel COM-interop script
# Applies market-price / profit value updates to the Ifrit_Profits leve-profit sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) as produced by the scheduled price-update runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 20
$ws.Range("H20").Value = 1640.3334
$ws.Range("I20").Value = 1640.3334
$ws.Range("K20").Value = 1640.3334
$ws.Range("M20").Value = -1410.3334
# ALC row 34
$ws.Range("H34").Value = 781.3333
$ws.Range("I34").Value = 781.3333
$ws.Range("K34").Value = 781.3333
$ws.Range("M34").Value = -578.3333
# ALC row 35
$ws.Range("H35").Value = 1640.3334
$ws.Range("I35").Value = 1640.3334
$ws.Range("K35").Value = 1640.3334
$ws.Range("M35").Value = -1261.3334
# ALC row 36
$ws.Range("H36").Value = 781.3333
$ws.Range("I36").Value = 781.3333
$ws.Range("K36").Value = 781.3333
$ws.Range("M36").Value = -66.33330000000001
# ALC row 40
$ws.Range("H40").Value = 1008.7059
$ws.Range("I40").Value = 977.2308
$ws.Range("K40").Value = 977.2308
$ws.Range("M40").Value = -802.2308

$ws = $wb.Worksheets.Item("ARM")
# ARM row 26
$ws.Range("H26").Value = 1416.3334
$ws.Range("I26").Value = 1416.3334
$ws.Range("K26").Value = 1416.3334
$ws.Range("M26").Value = -1086.3334
# ARM row 38
$ws.Range("H38").Value = 900.3333
$ws.Range("I38").Value = 900.3333
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 900.3333
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -433.3333
$ws.Range("N38").ClearContents()
# ARM row 39
$ws.Range("H39").Value = 5749.5
$ws.Range("I39").Value = 5749.5
$ws.Range("K39").Value = 5749.5
$ws.Range("M39").Value = -5229.5
# ARM row 61
$ws.Range("H61").Value = 2535.0527
$ws.Range("I61").Value = 1508.625
$ws.Range("K61").Value = 1508.625
$ws.Range("M61").Value = -1296.625
# ARM row 132
$ws.Range("H132").Value = 1280047.1
$ws.Range("I132").Value = 1570012.5
$ws.Range("K132").Value = 4710037.5
$ws.Range("M132").Value = -4707507.5
# ARM row 136
$ws.Range("H136").Value = 2535.0527
$ws.Range("I136").Value = 1508.625
$ws.Range("K136").Value = 4525.875
$ws.Range("M136").Value = -1975.875

$ws = $wb.Worksheets.Item("BSM")
# BSM row 107
$ws.Range("H107").Value = 221322.75
$ws.Range("I107").Value = 307364.1
$ws.Range("J107").Value = 1439.2778
$ws.Range("K107").Value = 307364.1
$ws.Range("L107").Value = 1439.2778
$ws.Range("M107").Value = -305444.1
$ws.Range("N107").Value = -5279.2778
# BSM row 134
$ws.Range("H134").Value = 45983.48
$ws.Range("I134").Value = 53653.305
$ws.Range("K134").Value = 160959.915
$ws.Range("M134").Value = -158424.915

$ws = $wb.Worksheets.Item("CRP")
# CRP row 58
$ws.Range("H58").Value = 3347.7026
$ws.Range("I58").Value = 1755.1333
$ws.Range("J58").Value = 4433.5454
$ws.Range("K58").Value = 1755.1333
$ws.Range("L58").Value = 4433.5454
$ws.Range("M58").Value = -1552.1333
$ws.Range("N58").Value = -4839.5454
# CRP row 92
$ws.Range("H92").Value = 27100.25
$ws.Range("J92").Value = 27100.25
$ws.Range("L92").Value = 27100.25
$ws.Range("N92").Value = -32092.25
# CRP row 134
$ws.Range("H134").Value = 2007.2444
$ws.Range("I134").Value = 2142.0789
$ws.Range("J134").Value = 1275.2858
$ws.Range("K134").Value = 6426.236699999999
$ws.Range("L134").Value = 3825.8574
$ws.Range("M134").Value = -3891.236699999999
$ws.Range("N134").Value = -8895.857400000001
# CRP row 135
$ws.Range("H135").Value = 40125
$ws.Range("J135").Value = 40125
$ws.Range("L135").Value = 40125
$ws.Range("N135").Value = -50265
# CRP row 136
$ws.Range("H136").Value = 3347.7026
$ws.Range("I136").Value = 1755.1333
$ws.Range("J136").Value = 4433.5454
$ws.Range("K136").Value = 5265.3999
$ws.Range("L136").Value = 13300.6362
$ws.Range("M136").Value = -2715.3999
$ws.Range("N136").Value = -18400.6362

$ws = $wb.Worksheets.Item("CUL")
# CUL row 49
$ws.Range("H49").Value = 2737.5
$ws.Range("J49").Value = 2737.5
$ws.Range("L49").Value = 8212.5
$ws.Range("N49").Value = -8524.5
# CUL row 102
$ws.Range("H102").Value = 4975
$ws.Range("J102").Value = 6300
$ws.Range("L102").Value = 18900
$ws.Range("N102").Value = -23768
# CUL row 103
$ws.Range("H103").Value = 8501575
$ws.Range("I103").Value = 11333767
$ws.Range("J103").Value = 5000
$ws.Range("K103").Value = 34001301
$ws.Range("L103").Value = 15000
$ws.Range("M103").Value = -34000422
$ws.Range("N103").Value = -16758
# CUL row 106
$ws.Range("H106").Value = 5910
$ws.Range("J106").Value = 5910
$ws.Range("L106").Value = 17730
$ws.Range("N106").Value = -19622
# CUL row 114
$ws.Range("H114").Value = 1843.75
$ws.Range("I114").Value = 863.3077
$ws.Range("J114").Value = 3002.4546
$ws.Range("K114").Value = 2589.9231
$ws.Range("L114").Value = 9007.363799999999
$ws.Range("M114").Value = 664.0769
$ws.Range("N114").Value = -15515.3638

$ws = $wb.Worksheets.Item("GSM")
# GSM row 92
$ws.Range("H92").Value = 6865
$ws.Range("J92").Value = 6865
$ws.Range("L92").Value = 6865
$ws.Range("N92").Value = -10609

$ws = $wb.Worksheets.Item("LTW")
# LTW row 22
$ws.Range("H22").Value = 290
$ws.Range("I22").Value = 266.66666
$ws.Range("J22").Value = 325
$ws.Range("K22").Value = 266.66666
$ws.Range("L22").Value = 325
$ws.Range("M22").Value = 28.33334000000002
$ws.Range("N22").Value = -915
# LTW row 27
$ws.Range("H27").Value = 290
$ws.Range("I27").Value = 266.66666
$ws.Range("J27").Value = 325
$ws.Range("K27").Value = 266.66666
$ws.Range("L27").Value = 325
$ws.Range("M27").Value = -159.66666
$ws.Range("N27").Value = -539
# LTW row 31
$ws.Range("H31").Value = 1359.1428
$ws.Range("I31").Value = 838
$ws.Range("J31").Value = 1750
$ws.Range("K31").Value = 838
$ws.Range("L31").Value = 1750
$ws.Range("M31").Value = -590
$ws.Range("N31").Value = -2246
# LTW row 76
$ws.Range("H76").Value = 8000
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 8000
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 8000
$ws.Range("N76").Value = -8676
$ws.Range("M76").ClearContents()
# LTW row 79
$ws.Range("H79").Value = 8000
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 8000
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 8000
$ws.Range("N79").Value = -10340
$ws.Range("M79").ClearContents()
# LTW row 122
$ws.Range("H122").Value = 6067.9688
$ws.Range("I122").Value = 6855.08
$ws.Range("K122").Value = 20565.24
$ws.Range("M122").Value = -18115.24
# LTW row 132
$ws.Range("H132").Value = 82368.28999999999
$ws.Range("I132").Value = 126039.78
$ws.Range("J132").Value = 3759.6
$ws.Range("K132").Value = 378119.34
$ws.Range("L132").Value = 11278.8
$ws.Range("M132").Value = -375589.34
$ws.Range("N132").Value = -16338.8
# LTW row 136
$ws.Range("H136").Value = 19766.666
$ws.Range("I136").Value = 34000
$ws.Range("J136").Value = 5533.3335
$ws.Range("K136").Value = 102000
$ws.Range("L136").Value = 16600.0005
$ws.Range("M136").Value = -99450
$ws.Range("N136").Value = -21700.0005

$ws = $wb.Worksheets.Item("WVR")
# WVR row 63
$ws.Range("H63").Value = 23448.75
$ws.Range("J63").Value = 29331.666
$ws.Range("L63").Value = 29331.666
$ws.Range("N63").Value = -30579.666
# WVR row 66
$ws.Range("H66").Value = 23448.75
$ws.Range("J66").Value = 29331.666
$ws.Range("L66").Value = 87994.99800000001
$ws.Range("N66").Value = -94234.99800000001
# WVR row 82
$ws.Range("H82").Value = 27650.5
$ws.Range("J82").Value = 27650.5
$ws.Range("L82").Value = 27650.5
$ws.Range("N82").Value = -28416.5
# WVR row 85
$ws.Range("H85").Value = 27650.5
$ws.Range("J85").Value = 27650.5
$ws.Range("L85").Value = 27650.5
$ws.Range("N85").Value = -30302.5
# WVR row 97
$ws.Range("H97").Value = 20285.5
$ws.Range("J97").Value = 20285.5
$ws.Range("L97").Value = 20285.5
$ws.Range("N97").Value = -22267.5
# WVR row 122
$ws.Range("H122").Value = 1974.9615
$ws.Range("I122").Value = 1368.5454
$ws.Range("J122").Value = 2419.6667
$ws.Range("K122").Value = 4105.6362
$ws.Range("L122").Value = 7259.000100000001
$ws.Range("M122").Value = -1655.6362
$ws.Range("N122").Value = -12159.0001
# WVR row 136
$ws.Range("H136").Value = 1156
$ws.Range("I136").Value = 972.8
$ws.Range("J136").Value = 1766.6666
$ws.Range("K136").Value = 2918.4
$ws.Range("L136").Value = 5299.9998
$ws.Range("M136").Value = -368.3999999999996
$ws.Range("N136").Value = -10399.9998

